$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6243
$ws.Range("I3").Value = 6516
$ws.Range("I4").Value = 1497
$ws.Range("I5").Value = 604
$ws.Range("I6").Value = 7382
$ws.Range("I7").Value = 22242

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range("I5").Value = 8
$ws.Range("I6").Value = 18

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I2").Value = 208
$ws.Range("I6").Value = 206
$ws.Range("I7").Value = 701

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I2").Value = 202
$ws.Range("I3").Value = 317
$ws.Range("I4").Value = 54
$ws.Range("I7").Value = 858

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I2").Value = 69
$ws.Range("I3").Value = 54
$ws.Range("I4").Value = 13
$ws.Range("I6").Value = 87
$ws.Range("I7").Value = 224

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 174
$ws.Range("I4").Value = 91
$ws.Range("I7").Value = 703
$ws.Range("I8").Value = 1338
$ws.Range("I11").Value = 332
$ws.Range("I15").Value = 255
$ws.Range("I18").Value = 166
$ws.Range("I19").Value = 621
$ws.Range("I20").Value = 558
$ws.Range("I23").Value = 220
$ws.Range("I24").Value = 63
$ws.Range("I29").Value = 1356
$ws.Range("I31").Value = 224
$ws.Range("I33").Value = 1010
$ws.Range("I35").Value = 30
$ws.Range("I36").Value = 301
$ws.Range("I37").Value = 701
$ws.Range("I38").Value = 18
$ws.Range("I41").Value = 95
$ws.Range("I42").Value = 782
$ws.Range("I44").Value = 164
$ws.Range("I47").Value = 161
$ws.Range("I48").Value = 295
$ws.Range("I51").Value = 264
$ws.Range("I52").Value = 479
$ws.Range("I53").Value = 240
$ws.Range("I54").Value = 458
$ws.Range("I60").Value = 124
$ws.Range("I67").Value = 858
$ws.Range("I69").Value = 49
$ws.Range("I72").Value = 87
$ws.Range("I73").Value = 204
$ws.Range("I76").Value = 319
$ws.Range("I78").Value = 302
$ws.Range("I79").Value = 627
$ws.Range("I80").Value = 73
$ws.Range("I82").Value = 26
$ws.Range("I83").Value = 484
$ws.Range("I85").Value = 1002
$ws.Range("I86").Value = 137
$ws.Range("I88").Value = 204
$ws.Range("I90").Value = 280
$ws.Range("I92").Value = 64
$ws.Range("I93").Value = 127
$ws.Range("I94").Value = 229
$ws.Range("I95").Value = 341
$ws.Range("I98").Value = 153
$ws.Range("I100").Value = 38
$ws.Range("I101").Value = 22242

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 165
$ws.Range("I7").Value = 484

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 118
$ws.Range("I7").Value = 341

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 229
$ws.Range("I3").Value = 376
$ws.Range("I4").Value = 44
$ws.Range("I7").Value = 1010

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I3").Value = 100
$ws.Range("I4").Value = 33
$ws.Range("I7").Value = 458

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I3").Value = 466
$ws.Range("I4").Value = 69
$ws.Range("I6").Value = 375
$ws.Range("I7").Value = 1356

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("I2").Value = 207
$ws.Range("I3").Value = 185
$ws.Range("I4").Value = 23
$ws.Range("I6").Value = 190
$ws.Range("I7").Value = 621

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I3").Value = 48
$ws.Range("I7").Value = 164

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I2").Value = 46
$ws.Range("I4").Value = 39
$ws.Range("I6").Value = 151
$ws.Range("I7").Value = 295

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I5").Value = 5
$ws.Range("I6").Value = 144
$ws.Range("I7").Value = 319

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I3").Value = 389
$ws.Range("I6").Value = 250
$ws.Range("I7").Value = 1002

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I2").Value = 29
$ws.Range("I7").Value = 95

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 196
$ws.Range("I3").Value = 242
$ws.Range("I6").Value = 265
$ws.Range("I7").Value = 782

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I3").Value = 77
$ws.Range("I7").Value = 302

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 63

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I6").Value = 65
$ws.Range("I7").Value = 220

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("I3").Value = 7
$ws.Range("I7").Value = 49

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I5").Value = 21
$ws.Range("I7").Value = 627

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 153
$ws.Range("I7").Value = 558

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I3").Value = 37
$ws.Range("I6").Value = 74
$ws.Range("I7").Value = 166

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I3").Value = 101
$ws.Range("I7").Value = 301

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("I2").Value = 35
$ws.Range("I3").Value = 32
$ws.Range("I7").Value = 127

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("I5").Value = 23
$ws.Range("I6").Value = 38

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I4").Value = 37
$ws.Range("I5").Value = 18
$ws.Range("I6").Value = 128
$ws.Range("I7").Value = 479

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I3").Value = 39
$ws.Range("I7").Value = 229

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I3").Value = 48
$ws.Range("I7").Value = 161

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I6").Value = 98
$ws.Range("I7").Value = 255

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 99
$ws.Range("I7").Value = 153

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I6").Value = 88
$ws.Range("I7").Value = 332

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("I3").Value = 9
$ws.Range("I7").Value = 30

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I2").Value = 67
$ws.Range("I3").Value = 64
$ws.Range("I7").Value = 204

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I3").Value = 56
$ws.Range("I7").Value = 174

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("I3").Value = 12
$ws.Range("I7").Value = 64

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I2").Value = 60
$ws.Range("I3").Value = 72
$ws.Range("I7").Value = 204

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I3").Value = 380
$ws.Range("I6").Value = 434
$ws.Range("I7").Value = 1338

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I2").Value = 26
$ws.Range("I4").Value = 63
$ws.Range("I7").Value = 137

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("I6").Value = 96
$ws.Range("I7").Value = 280

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I2").Value = 55
$ws.Range("I3").Value = 72
$ws.Range("I6").Value = 106
$ws.Range("I7").Value = 264

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I3").Value = 32
$ws.Range("I6").Value = 38
$ws.Range("I7").Value = 124

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I2").Value = 54
$ws.Range("I3").Value = 52
$ws.Range("I7").Value = 240

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I6").Value = 42
$ws.Range("I7").Value = 87

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("I5").Value = 15
$ws.Range("I6").Value = 26

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("I2").Value = 11
$ws.Range("I7").Value = 73

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 217
$ws.Range("I4").Value = 39
$ws.Range("I6").Value = 185
$ws.Range("I7").Value = 703

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("I3").Value = 24
$ws.Range("I7").Value = 91
